$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('C2').Value = '42 30 сер груз'
$ws.Range('G2').Value = '42, 30, сер, груз'
$ws.Range('I2').Value = 45341

# Row 3
$ws.Range('C3').Value = 'сер легк б/к'
$ws.Range('G3').Value = '42, 30, сер, груз'
$ws.Range('I3').Value = 45341

# Row 4
$ws.Range('C4').Value = 'сер легк б/к'
$ws.Range('G4').Value = 'сер, легк, б/к'
$ws.Range('I4').Value = 45341

# Row 5
$ws.Range('C5').Value = 'сер ошип'
$ws.Range('G5').Value = 'сер, легк, б/к'
$ws.Range('I5').Value = 45341

# Row 6
$ws.Range('C6').Value = 'сер легк'
$ws.Range('G6').Value = 'сер, ошип'
$ws.Range('I6').Value = 45341

# Row 7
$ws.Range('C7').Value = 'сер легк'
$ws.Range('G7').Value = 'сер, легк'
$ws.Range('I7').Value = 45341

# Row 8
$ws.Range('C8').Value = 'сер легк'
$ws.Range('G8').Value = 'сер, легк'
$ws.Range('I8').Value = 45341

# Row 9
$ws.Range('C9').Value = 'сер легк'
$ws.Range('G9').Value = 'сер, легк'
$ws.Range('I9').Value = 45341

# Row 10
$ws.Range('C10').Value = '210B Type сер H C'
$ws.Range('G10').Value = 'сер, легк'
$ws.Range('I10').Value = 45341

# Row 11
$ws.Range('C11').Value = 'Type сер LS-2 груз'
$ws.Range('G11').Value = '210B, Type, сер, H, C'
$ws.Range('I11').Value = 45341

# Row 12
$ws.Range('C12').Value = '202B Type сер C'
$ws.Range('G12').Value = '210B, Type, сер, H, C'
$ws.Range('I12').Value = 45341

# Row 13
$ws.Range('C13').Value = '202B Type сер LS-2 H C'
$ws.Range('G13').Value = 'Type, сер, LS-2, груз'
$ws.Range('I13').Value = 45341

# Row 14
$ws.Range('C14').Value = 'сер б/к груз'
$ws.Range('G14').Value = '202B, Type, сер, C'
$ws.Range('I14').Value = 45341

# Row 15
$ws.Range('C15').Value = 'сер легк б/к'
$ws.Range('G15').Value = '202B, Type, сер, LS-2, H, C'
$ws.Range('I15').Value = 45341

# Row 16
$ws.Range('C16').Value = 'сер легк б/к'
$ws.Range('G16').Value = '202B, Type, сер, LS-2, H, C'
$ws.Range('I16').Value = 45341

# Row 17
$ws.Range('C17').Value = '8 сер сх'
$ws.Range('G17').Value = '202B, Type, сер, LS-2, H, C'
$ws.Range('I17').Value = 45341

# Row 18
$ws.Range('C18').Value = 'сер легк'
$ws.Range('G18').Value = 'сер, б/к, груз'
$ws.Range('I18').Value = 45341

# Row 19
$ws.Range('C19').Value = 'сер легк'
$ws.Range('G19').Value = 'сер, б/к, груз'
$ws.Range('I19').Value = 45341

# Row 20
$ws.Range('C20').Value = 'сер легк'
$ws.Range('G20').Value = 'сер, б/к, груз'
$ws.Range('I20').Value = 45341

# Row 21
$ws.Range('C21').Value = 'сер легк'
$ws.Range('G21').Value = 'сер, б/к, груз'
$ws.Range('I21').Value = 45341

# Row 22
$ws.Range('C22').Value = 'сер легк'
$ws.Range('G22').Value = 'сер, легк, б/к'
$ws.Range('I22').Value = 45341

# Row 23
$ws.Range('G23').Value = 'сер, легк, б/к'
$ws.Range('I23').Value = 45341

# Row 24
$ws.Range('G24').Value = 'сер, легк'
$ws.Range('I24').Value = 45341

# Row 25
$ws.Range('G25').Value = 'сер, легк'
$ws.Range('I25').Value = 45341

# Row 26
$ws.Range('G26').Value = 'сер, легк'
$ws.Range('I26').Value = 45341

# Row 27
$ws.Range('G27').Value = 'сер, легк'
$ws.Range('I27').Value = 45341

# Row 28
$ws.Range('G28').Value = 'сер, легк'
$ws.Range('I28').Value = 45341

# Row 29
$ws.Range('G29').Value = 'сер, легк'
$ws.Range('I29').Value = 45341

# Row 30
$ws.Range('G30').Value = 'сер, легк'
$ws.Range('I30').Value = 45341
